$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 19309662
$ws.Range("I132").Value = 20491344
$ws.Range("K132").Value = 61474032
$ws.Range("M132").Value = -61471502
$ws.Range("H138").Value = 2357.0322
$ws.Range("I138").Value = 1307.5641
$ws.Range("J138").Value = 3114.9814
$ws.Range("K138").Value = 3922.6923
$ws.Range("L138").Value = 9344.9442
$ws.Range("M138").Value = 1217.3077
$ws.Range("N138").Value = -19624.9442

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 587.4737
$ws.Range("I2").Value = 550.7646999999999
$ws.Range("J2").Value = 899.5
$ws.Range("K2").Value = 550.7646999999999
$ws.Range("L2").Value = 899.5
$ws.Range("M2").Value = -437.7646999999999
$ws.Range("N2").Value = -1125.5
$ws.Range("H32").Value = 9611.819
$ws.Range("I32").Value = 6424.255
$ws.Range("J32").Value = 14692
$ws.Range("K32").Value = 6424.255
$ws.Range("L32").Value = 14692
$ws.Range("M32").Value = -6137.255
$ws.Range("N32").Value = -15266
$ws.Range("H45").Value = 1331.0667
$ws.Range("I45").Value = 833
$ws.Range("J45").Value = 2493.2222
$ws.Range("K45").Value = 833
$ws.Range("L45").Value = 2493.2222
$ws.Range("M45").Value = -456
$ws.Range("N45").Value = -3247.2222
$ws.Range("H61").Value = 2593
$ws.Range("I61").Value = 1980.2727
$ws.Range("K61").Value = 1980.2727
$ws.Range("M61").Value = -1768.2727
$ws.Range("H116").Value = 587.4737
$ws.Range("I116").Value = 550.7646999999999
$ws.Range("J116").Value = 899.5
$ws.Range("K116").Value = 550.7646999999999
$ws.Range("L116").Value = 899.5
$ws.Range("M116").Value = 1743.2353
$ws.Range("N116").Value = -5487.5
$ws.Range("H122").Value = 2481.5386
$ws.Range("I122").Value = 1340
$ws.Range("J122").Value = 5050
$ws.Range("K122").Value = 4020
$ws.Range("L122").Value = 15150
$ws.Range("M122").Value = -1570
$ws.Range("N122").Value = -20050
$ws.Range("H136").Value = 2593
$ws.Range("I136").Value = 1980.2727
$ws.Range("K136").Value = 5940.8181
$ws.Range("M136").Value = -3390.8181
$ws.Range("H137").Value = 50233.332
$ws.Range("J137").Value = 50233.332
$ws.Range("L137").Value = 50233.332
$ws.Range("N137").Value = -60433.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 587.4737
$ws.Range("I3").Value = 550.7646999999999
$ws.Range("J3").Value = 899.5
$ws.Range("K3").Value = 550.7646999999999
$ws.Range("L3").Value = 899.5
$ws.Range("M3").Value = -436.7646999999999
$ws.Range("N3").Value = -1127.5
$ws.Range("H86").Value = 2093.9048
$ws.Range("I86").Value = 2148.7144
$ws.Range("J86").Value = 1984.2858
$ws.Range("K86").Value = 2148.7144
$ws.Range("L86").Value = 1984.2858
$ws.Range("M86").Value = -1025.7144
$ws.Range("N86").Value = -4230.2858
$ws.Range("H89").Value = 2093.9048
$ws.Range("I89").Value = 2148.7144
$ws.Range("J89").Value = 1984.2858
$ws.Range("K89").Value = 10743.572
$ws.Range("L89").Value = 9921.429
$ws.Range("M89").Value = -5127.572
$ws.Range("N89").Value = -21153.429
$ws.Range("H105").Value = 2515.1365
$ws.Range("I105").Value = 2585.3684
$ws.Range("J105").Value = 2070.3333
$ws.Range("K105").Value = 2585.3684
$ws.Range("L105").Value = 2070.3333
$ws.Range("M105").Value = -838.3683999999998
$ws.Range("N105").Value = -5564.3333
$ws.Range("H134").Value = 2057.0833
$ws.Range("I134").Value = 1220.6528
$ws.Range("J134").Value = 7075.6665
$ws.Range("K134").Value = 3661.9584
$ws.Range("L134").Value = 21226.9995
$ws.Range("M134").Value = -1126.9584
$ws.Range("N134").Value = -26296.9995
$ws.Range("H137").Value = 32957.5
$ws.Range("J137").Value = 32957.5
$ws.Range("L137").Value = 32957.5
$ws.Range("N137").Value = -43157.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2593.451
$ws.Range("I31").Value = 988.9655
$ws.Range("J31").Value = 4708.4546
$ws.Range("K31").Value = 988.9655
$ws.Range("L31").Value = 4708.4546
$ws.Range("M31").Value = -693.9655
$ws.Range("N31").Value = -5298.4546
$ws.Range("H34").Value = 2593.451
$ws.Range("I34").Value = 988.9655
$ws.Range("J34").Value = 4708.4546
$ws.Range("K34").Value = 988.9655
$ws.Range("L34").Value = 4708.4546
$ws.Range("M34").Value = -786.9655
$ws.Range("N34").Value = -5112.4546
$ws.Range("H81").Value = 26000
$ws.Range("J81").Value = 26000
$ws.Range("L81").Value = 26000
$ws.Range("N81").Value = -27996
$ws.Range("H84").Value = 26000
$ws.Range("J84").Value = 26000
$ws.Range("L84").Value = 78000
$ws.Range("N84").Value = -87984

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 10870647
$ws.Range("I131").Value = 125002720
$ws.Range("J131").Value = 926.119
$ws.Range("K131").Value = 375008160
$ws.Range("L131").Value = 2778.357
$ws.Range("M131").Value = -375003120
$ws.Range("N131").Value = -12858.357

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5748.25
$ws.Range("I70").Value = 5397.278
$ws.Range("J70").Value = 6801.1665
$ws.Range("K70").Value = 5397.278
$ws.Range("L70").Value = 6801.1665
$ws.Range("M70").Value = -5127.278
$ws.Range("N70").Value = -7341.1665
$ws.Range("H73").Value = 5748.25
$ws.Range("I73").Value = 5397.278
$ws.Range("J73").Value = 6801.1665
$ws.Range("K73").Value = 5397.278
$ws.Range("L73").Value = 6801.1665
$ws.Range("M73").Value = -4461.278
$ws.Range("N73").Value = -8673.166499999999
$ws.Range("H102").Value = 2204.8696
$ws.Range("I102").Value = 1611.8334
$ws.Range("K102").Value = 1611.8334
$ws.Range("M102").Value = 10.16660000000002
$ws.Range("H122").Value = 2334.6743
$ws.Range("I122").Value = 1957.1714
$ws.Range("J122").Value = 3986.25
$ws.Range("K122").Value = 5871.5142
$ws.Range("L122").Value = 11958.75
$ws.Range("M122").Value = -3421.5142
$ws.Range("N122").Value = -16858.75
$ws.Range("H123").Value = 15505
$ws.Range("J123").Value = 15505
$ws.Range("L123").Value = 15505
$ws.Range("N123").Value = -20405
$ws.Range("H126").Value = 3821.8877
$ws.Range("I126").Value = 2748.2205
$ws.Range("J126").Value = 5446.154
$ws.Range("K126").Value = 8244.6615
$ws.Range("L126").Value = 16338.462
$ws.Range("M126").Value = -5774.6615
$ws.Range("N126").Value = -21278.462
$ws.Range("H132").Value = 2999.3103
$ws.Range("I132").Value = 1912.6842
$ws.Range("J132").Value = 5063.9
$ws.Range("K132").Value = 5738.0526
$ws.Range("L132").Value = 15191.7
$ws.Range("M132").Value = -3208.0526
$ws.Range("N132").Value = -20251.7
$ws.Range("H136").Value = 13160.52
$ws.Range("J136").Value = 13160.52
$ws.Range("L136").Value = 39481.56
$ws.Range("N136").Value = -44581.56

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5744.154
$ws.Range("J7").Value = 7413.5713
$ws.Range("L7").Value = 7413.5713
$ws.Range("N7").Value = -7637.5713
$ws.Range("H126").Value = 5744.154
$ws.Range("J126").Value = 7413.5713
$ws.Range("L126").Value = 22240.7139
$ws.Range("N126").Value = -27180.7139
$ws.Range("H132").Value = 5823.561
$ws.Range("I132").Value = 1753.2
$ws.Range("K132").Value = 5259.6
$ws.Range("M132").Value = -2729.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 29222080
$ws.Range("I81").Value = 32144086
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 64288172
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -64287111
$ws.Range("N81").Value = -6122
$ws.Range("H84").Value = 29222080
$ws.Range("I84").Value = 32144086
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 321440860
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -321435556
$ws.Range("N84").Value = -30608
$ws.Range("H113").Value = 277.67856
$ws.Range("I113").Value = 261.05554
$ws.Range("K113").Value = 783.16662
$ws.Range("M113").Value = 1386.83338
$ws.Range("H122").Value = 3129.3076
$ws.Range("I122").Value = 2353.9614
$ws.Range("K122").Value = 7061.8842
$ws.Range("M122").Value = -4611.8842
$ws.Range("H126").Value = 293070.22
$ws.Range("I126").Value = 2398.9546
$ws.Range("J126").Value = 719388.0600000001
$ws.Range("K126").Value = 7196.8638
$ws.Range("L126").Value = 2158164.18
$ws.Range("M126").Value = -4726.8638
$ws.Range("N126").Value = -2163104.18
$ws.Range("H132").Value = 6174464
$ws.Range("I132").Value = 583.26086
$ws.Range("J132").Value = 10755085
$ws.Range("K132").Value = 1749.78258
$ws.Range("L132").Value = 32265255
$ws.Range("M132").Value = 780.2174199999999
$ws.Range("N132").Value = -32270315
$ws.Range("H136").Value = 2625.932
$ws.Range("I136").Value = 742.26666
$ws.Range("J136").Value = 6662.357
$ws.Range("K136").Value = 2226.79998
$ws.Range("L136").Value = 19987.071
$ws.Range("M136").Value = 323.2000200000002
$ws.Range("N136").Value = -25087.071
